# Update dSF column (F) values to reflect re-pulled source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -5
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 6
$ws.Range("F32").Value = -7
$ws.Range("F37").Value = -3
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = -9
$ws.Range("F40").Value = 4
$ws.Range("F41").Value = 2
